$d = $word.ActiveDocument

$replacements = @(
    @("318×8=", "472×6="),
    @("270×2=", "678×9="),
    @("711×5=", "817×5="),
    @("254×9=", "860×7="),
    @("887×3=", "366×6="),
    @("816×5=", "913×7="),
    @("514×8=", "256×8="),
    @("708×9=", "717×9="),
    @("842×8=", "584×4="),
    @("604×2=", "226×9="),
    @("536×9=", "728×4="),
    @("189×5=", "648×5="),
    @("957×8=", "744×3="),
    @("301×2=", "612×7="),
    @("664×2=", "947×3="),
    @("370×3=", "742×6="),
    @("288×3=", "275×9="),
    @("580×3=", "234×8="),
    @("651×8=", "727×4="),
    @("273×5=", "666×4="),
    @("794×9=", "283×7="),
    @("859×2=", "195×2="),
    @("158×2=", "879×5="),
    @("960×6=", "213×5="),
    @("939×8=", "299×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
